{"js": "const replacements = [\n  [\"657\u00d73=\", \"892\u00d77=\"],\n  [\"227\u00d79=\", \"319\u00d72=\"],\n  [\"116\u00d77=\", \"184\u00d77=\"],\n  [\"782\u00d79=\", \"189\u00d72=\"],\n  [\"852\u00d74=\", \"321\u00d76=\"],\n  [\"328\u00d79=\", \"978\u00d76=\"],\n  [\"591\u00d77=\", \"107\u00d76=\"],\n  [\"153\u00d74=\", \"866\u00d77=\"],\n  [\"529\u00d72=\", \"686\u00d76=\"],\n  [\"846\u00d74=\", \"264\u00d79=\"],\n  [\"149\u00d76=\", \"162\u00d73=\"],\n  [\"196\u00d74=\", \"472\u00d75=\"],\n  [\"522\u00d77=\", \"340\u00d75=\"],\n  [\"217\u00d78=\", \"101\u00d79=\"],\n  [\"286\u00d77=\", \"114\u00d78=\"],\n  [\"536\u00d76=\", \"793\u00d77=\"],\n  [\"906\u00d78=\", \"161\u00d73=\"],\n  [\"906\u00d76=\", \"520\u00d72=\"],\n  [\"599\u00d72=\", \"684\u00d72=\"],\n  [\"493\u00d77=\", \"593\u00d79=\"],\n  [\"117\u00d76=\", \"296\u00d75=\"],\n  [\"870\u00d72=\", \"859\u00d72=\"],\n  [\"247\u00d72=\", \"172\u00d73=\"],\n  [\"976\u00d75=\", \"566\u00d73=\"],\n  [\"714\u00d74=\", \"243\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old=\"657\u00d73=\"; new=\"892\u00d77=\"},\n    @{old=\"227\u00d79=\"; new=\"319\u00d72=\"},\n    @{old=\"116\u00d77=\"; new=\"184\u00d77=\"},\n    @{old=\"782\u00d79=\"; new=\"189\u00d72=\"},\n    @{old=\"852\u00d74=\"; new=\"321\u00d76=\"},\n    @{old=\"328\u00d79=\"; new=\"978\u00d76=\"},\n    @{old=\"591\u00d77=\"; new=\"107\u00d76=\"},\n    @{old=\"153\u00d74=\"; new=\"866\u00d77=\"},\n    @{old=\"529\u00d72=\"; new=\"686\u00d76=\"},\n    @{old=\"846\u00d74=\"; new=\"264\u00d79=\"},\n    @{old=\"149\u00d76=\"; new=\"162\u00d73=\"},\n    @{old=\"196\u00d74=\"; new=\"472\u00d75=\"},\n    @{old=\"522\u00d77=\"; new=\"340\u00d75=\"},\n    @{old=\"217\u00d78=\"; new=\"101\u00d79=\"},\n    @{old=\"286\u00d77=\"; new=\"114\u00d78=\"},\n    @{old=\"536\u00d76=\"; new=\"793\u00d77=\"},\n    @{old=\"906\u00d78=\"; new=\"161\u00d73=\"},\n    @{old=\"906\u00d76=\"; new=\"520\u00d72=\"},\n    @{old=\"599\u00d72=\"; new=\"684\u00d72=\"},\n    @{old=\"493\u00d77=\"; new=\"593\u00d79=\"},\n    @{old=\"117\u00d76=\"; new=\"296\u00d75=\"},\n    @{old=\"870\u00d72=\"; new=\"859\u00d72=\"},\n    @{old=\"247\u00d72=\"; new=\"172\u00d73=\"},\n    @{old=\"976\u00d75=\"; new=\"566\u00d73=\"},\n    @{old=\"714\u00d74=\"; new=\"243\u00d74=\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.new\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $r.new, 2) | Out-Null\n}\n"}
